# "Fixed typos for 2019"
# On slide 5, the sentence under the pandas/table diagram reads
# "...we use pandas to retrieve..." - capitalise "pandas" -> "Pandas"
# (the library's proper name), matching the author's correction.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item("TextBox 11")
$tr = $sh.TextFrame.TextRange

# Re-type the middle portion of the sentence ("use pandas " -> "use Pandas ")
# by selecting exactly that substring and overwriting it - this is how the
# author would have fixed the capitalisation by hand, and it naturally
# splits the single run into the three runs PowerPoint keeps around the
# edited word.
$fragment = $tr.Characters(20, 11)
$fragment.Text = "use Pandas "
